$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: num_customers 59 -> 60, retention_rate recalculated (C27/D27)
$ws.Range("C27").Value = 60
$ws.Range("E27").Value = 0.02664298401420959

# Row 36: num_customers 152 -> 153, retention_rate recalculated (C36/D36)
$ws.Range("C36").Value = 153
$ws.Range("E36").Value = 0.07927461139896373

# Row 37: num_customers 1025 -> 1028, cohort_size 1025 -> 1028 (E37 stays 1)
$ws.Range("C37").Value = 1028
$ws.Range("D37").Value = 1028
